# Mise à jour de l'application
# Add 13 new training-log rows (245-257) for the 2025-09-02 session, extend the
# charge ("I") formula down to the new rows, and refresh the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the current last row (244).
# Columns: Date(serial), Name, Volume, Intensite, Fatigue, Douleur, Localisation, Plaisir
$newRows = @(
    @(45902, "Amir Etien",       70, 5, 6, 0, "",         4),
    @(45902, "Omar Benyounes",   70, 6, 6, 0, "",         6),
    @(45902, "Rayane Chayebi",   70, 6, 6, 6, "Adducteur", 6),
    @(45902, "Yoan Zouma",       70, 3, 6, 3, "Cheville",  6),
    @(45902, "Ilyes Boughanmi",  70, 6, 5, 0, "",         0),
    @(45902, "Jeremie Laurent",  70, 5, 5, 0, "",         8),
    @(45902, "Amir Kherrab",     70, 5, 5, 5, "Semelle ",  7),
    @(45902, "Naim Ighbane",     70, 6, 0, 0, "",         9),
    @(45902, "Emmanuel Valey",   70, 3, 2, 4, "Adducteur", 0),
    @(45902, "Ilan Ihaddadene",  75, 5, 5, 2, "Semelle",   7),
    @(45902, "Karahali Souaré",  70, 2, 5, 7, "Cheville",  8),
    @(45902, "Sofiane Belle",    70, 5, 3, 0, "",         6),
    @(45902, "Mattheo Haon",     70, 6, 7, 5, "Adducteur ", 6)
)

$firstNewRow = 245
$lastExisting = 244
$lastNewRow = $firstNewRow + $newRows.Count - 1

$emptyLocRow = 241  # known row whose "Localisation douleur" cell is blank (centered style)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $data = $newRows[$i]

    # Duplicate formatting (fonts / number formats / alignment) from the last
    # existing data row so the new row matches the table's look.
    $ws.Range("A" + $lastExisting + ":I" + $lastExisting).Copy($ws.Range("A" + $r + ":I" + $r))

    $ws.Range("A" + $r).Value = $data[0]
    $ws.Range("B" + $r).Value = $data[1]
    $ws.Range("C" + $r).Value = $data[2]
    $ws.Range("D" + $r).Value = $data[3]
    $ws.Range("E" + $r).Value = $data[4]
    $ws.Range("F" + $r).Value = $data[5]

    if ([string]::IsNullOrEmpty($data[6])) {
        $ws.Range("G" + $r).ClearContents()
        # Re-apply the blank cell's centered style (copy wipes it to the text style).
        $ws.Range("G" + $emptyLocRow).Copy($ws.Range("G" + $r))
    } else {
        $ws.Range("G" + $r).Value = $data[6]
    }

    $ws.Range("H" + $r).Value = $data[7]
    $ws.Range("I" + $r).Formula = "=C" + $r + "*D" + $r
}

# Scroll / select near the newly-added rows, mirroring where the author
# was working when the rows were added.
$excel.ActiveWindow.ScrollRow = 230
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L252").Select()
